# Add files via upload
#
# The workbook originally had a single sheet "rules_flat" holding a flat
# allocation-rules table (one parent OPK -> one child OPK per row). This
# change adds two new source sheets in front of it:
#   - parent_header : the parent OPK (cost centre) header/summary row
#   - children_table: the child OPK rows that get expanded into rules_flat
# and rebuilds "rules_flat" (still the last tab) from those children, one
# row per child, with the weight coming from the children_table KVI column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New sheet 1: parent_header (inserted before rules_flat)
# ---------------------------------------------------------------------
$parentHeader = $wb.Worksheets.Add($null, $null, 1, $null)
$parentHeader.Name = "parent_header"

$parentHeader.Range("A1").Value = "Id OPK (rodzic)"
$parentHeader.Range("B1").Value = "OPK (nazwa)"
$parentHeader.Range("C1").Value = "Kwota"
$parentHeader.Range("A1:C1").Font.Bold = $true
$parentHeader.Range("A1:C1").HorizontalAlignment = -4108
$parentHeader.Range("A1:C1").VerticalAlignment = -4160
$parentHeader.Range("A1:C1").Borders.LineStyle = 1

# A2 ("73") must stay text, not be auto-converted to a number.
$parentHeader.Range("A2").NumberFormat = "@"
$parentHeader.Range("A2").Value = "73"
$parentHeader.Range("A2").Style = "Normal"
$parentHeader.Range("B2").Value = "Wydział MO"
$parentHeader.Range("C2").Value = 100000

# ---------------------------------------------------------------------
# New sheet 2: children_table (inserted between parent_header and rules_flat)
# ---------------------------------------------------------------------
$childrenTable = $wb.Worksheets.Add($null, $parentHeader)
$childrenTable.Name = "children_table"

$childrenTable.Range("A1").Value = "Id OPK (dziecko)"
$childrenTable.Range("B1").Value = "OPK (nazwa)"
$childrenTable.Range("C1").Value = "Kwota"
$childrenTable.Range("D1").Value = "KVI"
$childrenTable.Range("A1:D1").Font.Bold = $true
$childrenTable.Range("A1:D1").HorizontalAlignment = -4108
$childrenTable.Range("A1:D1").VerticalAlignment = -4160
$childrenTable.Range("A1:D1").Borders.LineStyle = 1

$childIds = @("73.54", "73.55", "73.516", "73.515")
$childNames = @("Stanowisko 54", "Stanowisko 55", "Stanowisko 516", "Stanowisko 515")
$childKvi = @(40, 30, 20, 10)

# Column A holds "73.54" style ids -- keep them text so they don't get
# mangled into floating point numbers (73.540000000000006 etc.).
$childrenTable.Range("A2:A5").NumberFormat = "@"
for ($i = 0; $i -lt 4; $i++) {
    $row = $i + 2
    $childrenTable.Cells.Item($row, 1).Value = $childIds[$i]
    $childrenTable.Cells.Item($row, 2).Value = $childNames[$i]
    $childrenTable.Cells.Item($row, 4).Value = $childKvi[$i]
}
$childrenTable.Range("A2:A5").Style = "Normal"

# ---------------------------------------------------------------------
# rules_flat: rebuild from the children, one row per child
# (re-fetched by name now, after the inserts, so it points at the right
#  physical sheet instead of a stale pre-insert position)
# ---------------------------------------------------------------------
$rulesFlat = $wb.Worksheets.Item("rules_flat")
$rulesFlat.Range("A2:G5").ClearContents()

$rulesFlat.Range("A2:B5").NumberFormat = "@"
for ($i = 0; $i -lt 4; $i++) {
    $row = $i + 2
    $rulesFlat.Cells.Item($row, 1).Value = "73"
    $rulesFlat.Cells.Item($row, 2).Value = $childIds[$i]
    $rulesFlat.Cells.Item($row, 3).Value = "KVI"
    $rulesFlat.Cells.Item($row, 4).Value = $childKvi[$i]
}
$rulesFlat.Range("A2:B5").Style = "Normal"

# The original sheet's tab-selected flag carries over to the physical
# sheet1.xml part, which is now "parent_header" (the new first tab).
$parentHeader.Activate()
